# Apply updated crypto price/volume snapshot to columns D (Price) and E (Volume 1h)
# for worksheet rows 2-51. Values are plain text in the source data (coinranking scrape),
# so numeric-looking prices are entered with a leading apostrophe to force text storage,
# then restyled back to "Normal" so no stray number-format style is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '51.037.50'
$ws.Range("E2").Value = '  -1.30%  '
# Row 3
$ws.Range("D3").Value = '2.751.30'
$ws.Range("E3").Value = '  -0.69%  '
# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '
# Row 5
$ws.Range("D5").Value = "'352.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.54%  '
# Row 6
$ws.Range("D6").Value = "'107.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.73%  '
# Row 7
$ws.Range("E7").Value = '  -2.31%  '
# Row 8
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.06%  '
# Row 9
$ws.Range("E9").Value = '  -1.91%  '
# Row 10
$ws.Range("D10").Value = "'39.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.97%  '
# Row 11
$ws.Range("E11").Value = '  +3.51%  '
# Row 12
$ws.Range("D12").Value = "'0.0829"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.07%  '
# Row 13
$ws.Range("D13").Value = "'19.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.81%  '
# Row 14
$ws.Range("E14").Value = '  -2.04%  '
# Row 15
$ws.Range("D15").Value = '3.174.61'
$ws.Range("E15").Value = '  -1.11%  '
# Row 16
$ws.Range("D16").Value = '2.740.23'
$ws.Range("E16").Value = '  -2.57%  '
# Row 17
$ws.Range("D17").Value = "'0.920"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.89%  '
# Row 18
$ws.Range("D18").Value = '50.979.23'
$ws.Range("E18").Value = '  -1.13%  '
# Row 19
$ws.Range("D19").Value = "'7.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.21%  '
# Row 20
$ws.Range("E20").Value = '  -2.41%  '
# Row 21
$ws.Range("D21").Value = "'12.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.93%  '
# Row 22
$ws.Range("E22").Value = '  -2.39%  '
# Row 23
$ws.Range("D23").Value = "'69.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.50%  '
# Row 24
$ws.Range("D24").Value = "'263.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.67%  '
# Row 25
$ws.Range("E25").Value = '  -1.58%  '
# Row 26
$ws.Range("E26").Value = '  +0.03%  '
# Row 27
$ws.Range("D27").Value = "'25.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.28%  '
# Row 28
$ws.Range("D28").Value = "'0.160"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +13.77%  '
# Row 29
$ws.Range("E29").Value = '  +0.22%  '
# Row 30
$ws.Range("D30").Value = "'10.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.16%  '
# Row 31
$ws.Range("D31").Value = "'51.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.60%  '
# Row 32
$ws.Range("D32").Value = "'34.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.66%  '
# Row 33
$ws.Range("D33").Value = "'5.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.44%  '
# Row 34
$ws.Range("D34").Value = "'0.0438"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.05%  '
# Row 35
$ws.Range("E35").Value = '  -1.42%  '
# Row 36
$ws.Range("E36").Value = '  -2.97%  '
# Row 37
$ws.Range("E37").Value = '  -0.22%  '
# Row 38
$ws.Range("D38").Value = "'18.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.93%  '
# Row 39
$ws.Range("E39").Value = '  -1.57%  '
# Row 40
$ws.Range("E40").Value = '  -2.75%  '
# Row 41
$ws.Range("E41").Value = '  -1.39%  '
# Row 42
$ws.Range("D42").Value = "'2.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.00%  '
# Row 43
$ws.Range("D43").Value = "'120.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.32%  '
# Row 44
$ws.Range("E44").Value = '  -2.30%  '
# Row 45
$ws.Range("D45").Value = "'21.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.47%  '
# Row 46
$ws.Range("D46").Value = '2.082.39'
$ws.Range("E46").Value = '  +1.34%  '
# Row 48
$ws.Range("D48").Value = "'2.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.30%  '
# Row 49
$ws.Range("D49").Value = "'0.911"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.57%  '
# Row 50
$ws.Range("E50").Value = '  -4.27%  '
# Row 51
$ws.Range("E51").Value = '  +5.07%  '
